$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "Match Date" column
$ws.Range("K1").Value = "Match Date"

# Column width for K (best-fit to the date values, mirrors Excel's AutoFit result)
$ws.Columns.Item(11).ColumnWidth = 8.63

# Date values (stored as serial numbers, formatted as plain integers -- numFmtId 1 "0")
# K5 is left blank (row 5 has no match data) but still picks up the column's number format.
$ws.Range("K3").Value = 42461
$ws.Range("K4").Value = 42464
$ws.Range("K6").Value = 42459
$ws.Range("K7").Value = 42467

$ws.Range("K3:K7").NumberFormat = "0"

# Update selection to match target state
$ws.Range("M4").Select()
